$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 for "Hampton" (shifts existing rows down by one)
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Hampton"
$ws.Range("B9").Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Range("C9").Value = "28-12-2020 1:20pm-2:30pm"
$ws.Range("D9").Value = "Case ate in store"

# Insert a new row at row 20 for "McKinnon" (shifts existing rows down by one)
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "McKinnon"
$ws.Range("B20").Value = "260 McKinnon Road, McKinnon VIC 3204"
$ws.Range("C20").Value = "23-12-2020 4:00pm-6:00pm"
$ws.Range("D20").Value = "Case had hair cut in store"
